# Update the date line in the title paragraph.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-05-16 Thursday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-05-17 Friday", 2) | Out-Null

# Update the division-problem answers in the single table. Each data row
# (1, 5, 9, 13, 17) holds five answer cells; addressing cells by
# (row, column) avoids any ambiguity from values that are reused as both
# an old and a new answer elsewhere in the table.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "63÷4=15, 3"
$t.Cell(1, 2).Range.Text = "70÷4=17, 2"
$t.Cell(1, 3).Range.Text = "95÷3=31, 2"
$t.Cell(1, 4).Range.Text = "19÷2=9, 1"
$t.Cell(1, 5).Range.Text = "65÷3=21, 2"

$t.Cell(5, 1).Range.Text = "28÷6=4, 4"
$t.Cell(5, 2).Range.Text = "99÷9=11, 0"
$t.Cell(5, 3).Range.Text = "97÷6=16, 1"
$t.Cell(5, 4).Range.Text = "74÷6=12, 2"
$t.Cell(5, 5).Range.Text = "74÷7=10, 4"

$t.Cell(9, 1).Range.Text = "21÷8=2, 5"
$t.Cell(9, 2).Range.Text = "79÷2=39, 1"
$t.Cell(9, 3).Range.Text = "57÷8=7, 1"
$t.Cell(9, 4).Range.Text = "13÷9=1, 4"
$t.Cell(9, 5).Range.Text = "37÷7=5, 2"

$t.Cell(13, 1).Range.Text = "57÷5=11, 2"
$t.Cell(13, 2).Range.Text = "15÷7=2, 1"
$t.Cell(13, 3).Range.Text = "92÷7=13, 1"
$t.Cell(13, 4).Range.Text = "72÷8=9, 0"
$t.Cell(13, 5).Range.Text = "79÷7=11, 2"

$t.Cell(17, 1).Range.Text = "35÷2=17, 1"
$t.Cell(17, 2).Range.Text = "88÷2=44, 0"
$t.Cell(17, 3).Range.Text = "71÷2=35, 1"
$t.Cell(17, 4).Range.Text = "95÷4=23, 3"
$t.Cell(17, 5).Range.Text = "44÷8=5, 4"
